$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 15790.143
$ws.Range("I9").Value = 25407.75
$ws.Range("M9").Value = -25238.75
$ws.Range("K9").Value = 25407.75
$ws.Range("I15").Value = 685.25
$ws.Range("M15").Value = -1886.75
$ws.Range("K15").Value = 2055.75
$ws.Range("H15").Value = 685.25
$ws.Range("H92").Value = 515.46155
$ws.Range("K92").Value = 471.63635
$ws.Range("M92").Value = 776.36365
$ws.Range("I92").Value = 471.63635
$ws.Range("N97").Value = -6506.3999
$ws.Range("J97").Value = 1838.1333
$ws.Range("H97").Value = 1838.1333
$ws.Range("L97").Value = 5514.3999
$ws.Range("L112").Value = 16860172.5
$ws.Range("J112").Value = 5620057.5
$ws.Range("N112").Value = -16862388.5
$ws.Range("H112").Value = 5557712
$ws.Range("M132").Value = -2071.634
$ws.Range("K132").Value = 4601.634
$ws.Range("I132").Value = 1533.878
$ws.Range("H132").Value = 2099.2246
$ws.Range("M137").Value = -3042.5358
$ws.Range("I137").Value = 1864.1786
$ws.Range("K137").Value = 5592.5358
$ws.Range("H137").Value = 1809.1111
$ws.Range("N138").Value = -76966058
$ws.Range("H138").Value = 19616650
$ws.Range("L138").Value = 76955778
$ws.Range("J138").Value = 25651926
$ws.Range("H141").Value = 1828.1818
$ws.Range("J141").Value = 3105
$ws.Range("I141").Value = 1700.5
$ws.Range("M141").Value = 78.5
$ws.Range("K141").Value = 5101.5
$ws.Range("N141").Value = -19675
$ws.Range("L141").Value = 9315

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 6204.7144
$ws.Range("M32").Value = -5917.7144
$ws.Range("H32").Value = 6542.483
$ws.Range("K32").Value = 6204.7144
$ws.Range("H43").Value = 200039200
$ws.Range("N43").Value = -49619.5
$ws.Range("L43").Value = 48993.5
$ws.Range("J43").Value = 48993.5
$ws.Range("H61").Value = 11115378
$ws.Range("M61").Value = -11908526
$ws.Range("K61").Value = 11908738
$ws.Range("I61").Value = 11908738
$ws.Range("H74").Value = 30337660
$ws.Range("M74").Value = -32293990
$ws.Range("I74").Value = 32294864
$ws.Range("K74").Value = 32294864
$ws.Range("K77").Value = 161474320
$ws.Range("M77").Value = -161469952
$ws.Range("I77").Value = 32294864
$ws.Range("H77").Value = 30337660
$ws.Range("K102").Value = 669329.7
$ws.Range("I102").Value = 669329.7
$ws.Range("M102").Value = -667707.7
$ws.Range("H102").Value = 289882.16
$ws.Range("L102").Value = 5296.5
$ws.Range("J102").Value = 5296.5
$ws.Range("N102").Value = -8540.5
$ws.Range("I122").Value = 2097.7727
$ws.Range("H122").Value = 2816.775
$ws.Range("M122").Value = -3843.3181
$ws.Range("K122").Value = 6293.3181
$ws.Range("H136").Value = 11115378
$ws.Range("M136").Value = -35723664
$ws.Range("K136").Value = 35726214
$ws.Range("I136").Value = 11908738
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("M141").Value = $null
$ws.Range("K141").Value = 0
$ws.Range("N141").Value = $null
$ws.Range("L141").Value = 0

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K31").Value = 2170.3901
$ws.Range("J31").Value = 8270.647000000001
$ws.Range("H31").Value = 3958.3965
$ws.Range("M31").Value = -1875.3901
$ws.Range("N31").Value = -8860.647000000001
$ws.Range("L31").Value = 8270.647000000001
$ws.Range("I31").Value = 2170.3901
$ws.Range("N34").Value = -8674.647000000001
$ws.Range("L34").Value = 8270.647000000001
$ws.Range("M34").Value = -1968.3901
$ws.Range("H34").Value = 3958.3965
$ws.Range("K34").Value = 2170.3901
$ws.Range("I34").Value = 2170.3901
$ws.Range("J34").Value = 8270.647000000001
$ws.Range("L58").Value = 3011.2222
$ws.Range("N58").Value = -3417.2222
$ws.Range("J58").Value = 3011.2222
$ws.Range("H58").Value = 1850
$ws.Range("K134").Value = 4116.620699999999
$ws.Range("H134").Value = 1368.0555
$ws.Range("N134").Value = -9122.571599999999
$ws.Range("L134").Value = 4052.5716
$ws.Range("J134").Value = 1350.8572
$ws.Range("I134").Value = 1372.2069
$ws.Range("M134").Value = -1581.620699999999
$ws.Range("J136").Value = 3011.2222
$ws.Range("L136").Value = 9033.6666
$ws.Range("H136").Value = 1850
$ws.Range("N136").Value = -14133.6666

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M8").Value = -1006.25
$ws.Range("I8").Value = 381.75
$ws.Range("H8").Value = 381.75
$ws.Range("K8").Value = 1145.25
$ws.Range("K134").Value = 8631
$ws.Range("H134").Value = 9881.684999999999
$ws.Range("N134").Value = -63133.99800000001
$ws.Range("L134").Value = 52993.99800000001
$ws.Range("J134").Value = 17664.666
$ws.Range("I134").Value = 2877
$ws.Range("M134").Value = -3561

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N104").Value = -55488
$ws.Range("L104").Value = 48500
$ws.Range("H104").Value = 48500
$ws.Range("J104").Value = 48500
$ws.Range("N132").Value = -13258.4999
$ws.Range("M132").Value = -9183.6158
$ws.Range("K132").Value = 11713.6158
$ws.Range("I132").Value = 3904.5386
$ws.Range("J132").Value = 2732.8333
$ws.Range("L132").Value = 8198.499899999999
$ws.Range("H132").Value = 3534.5264

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K46").Value = 922.7778
$ws.Range("M46").Value = -734.7778
$ws.Range("H46").Value = 1792.3334
$ws.Range("I46").Value = 922.7778
$ws.Range("L46").Value = 4401
$ws.Range("N46").Value = -4777
$ws.Range("J46").Value = 4401
$ws.Range("H61").Value = 2148.8936
$ws.Range("M61").Value = -1300.1666
$ws.Range("L61").Value = 3290.1765
$ws.Range("N61").Value = -3694.1765
$ws.Range("K61").Value = 1502.1666
$ws.Range("I61").Value = 1502.1666
$ws.Range("J61").Value = 3290.1765
$ws.Range("I82").Value = 0
$ws.Range("M82").Value = $null
$ws.Range("H82").Value = 4750
$ws.Range("N82").Value = -5472
$ws.Range("L82").Value = 4750
$ws.Range("J82").Value = 4750
$ws.Range("K82").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("N85").Value = -7246
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 4750
$ws.Range("L85").Value = 4750
$ws.Range("M85").Value = $null
$ws.Range("H85").Value = 4750
$ws.Range("H105").Value = 73166
$ws.Range("J105").Value = 73166
$ws.Range("N105").Value = -80154
$ws.Range("L105").Value = 73166
$ws.Range("K113").Value = 1502.1666
$ws.Range("M113").Value = 667.8334
$ws.Range("H113").Value = 2148.8936
$ws.Range("L113").Value = 3290.1765
$ws.Range("I113").Value = 1502.1666
$ws.Range("J113").Value = 3290.1765
$ws.Range("N113").Value = -7630.1765
$ws.Range("I122").Value = 2717.7693
$ws.Range("J122").Value = 5158.909
$ws.Range("H122").Value = 4252.2
$ws.Range("M122").Value = -5703.3079
$ws.Range("N122").Value = -20376.727
$ws.Range("L122").Value = 15476.727
$ws.Range("K122").Value = 8153.3079

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 26052.059
$ws.Range("J51").Value = 29332.5
$ws.Range("M51").Value = -17669
$ws.Range("L51").Value = 29332.5
$ws.Range("I51").Value = 18179
$ws.Range("N51").Value = -30352.5
$ws.Range("K51").Value = 18179
$ws.Range("H52").Value = 34507.05
$ws.Range("M52").Value = -13686.25
$ws.Range("K52").Value = 13912.25
$ws.Range("I52").Value = 13912.25
$ws.Range("J62").Value = 7210.5557
$ws.Range("M62").Value = -2360
$ws.Range("K62").Value = 2984
$ws.Range("L62").Value = 7210.5557
$ws.Range("H62").Value = 5701.0713
$ws.Range("N62").Value = -8458.555700000001
$ws.Range("I62").Value = 2984
$ws.Range("I65").Value = 2984
$ws.Range("L65").Value = 36052.7785
$ws.Range("J65").Value = 7210.5557
$ws.Range("N65").Value = -42292.7785
$ws.Range("H65").Value = 5701.0713
$ws.Range("K65").Value = 14920
$ws.Range("M65").Value = -11800
$ws.Range("I107").Value = 463.57144
$ws.Range("H107").Value = 449.44446
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 1390.71432
$ws.Range("M107").Value = 529.28568
$ws.Range("L107").Value = 1200
$ws.Range("N107").Value = -5040
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
$ws.Range("N132").Value = -8045
$ws.Range("M132").Value = -7066.499899999999
$ws.Range("K132").Value = 9596.499899999999
$ws.Range("I132").Value = 3198.8333
$ws.Range("J132").Value = 995
$ws.Range("L132").Value = 2985
$ws.Range("H132").Value = 3061.0938
$ws.Range("J136").Value = 3200
$ws.Range("L136").Value = 9600
$ws.Range("H136").Value = 930.5161000000001
$ws.Range("M136").Value = -14.59991999999966
$ws.Range("N136").Value = -14700
$ws.Range("K136").Value = 2564.59992
$ws.Range("I136").Value = 854.86664

Write-Output "Applied all changes"